# Applies numeric corrections to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 25 - West Ham United
$ws.Range("I25").Value = 59274

# Row 28 - Leicester City
$ws.Range("R28").Value = 28
$ws.Range("W28").Value = 0

# Row 33 - Wolverhampton Wanderers
$ws.Range("I33").Value = 86542

# Row 34 - Olympique Marseille
$ws.Range("I34").Value = 76530

# Row 35 - Real Sociedad
$ws.Range("I35").Value = 80396

# Row 36 - Eintracht Frankfurt
$ws.Range("U36").Value = 18

# Row 39 - Nice
$ws.Range("C39").Value = 9
$ws.Range("D39").Value = 9
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 3.48
$ws.Range("H39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("U39").Value = 8

# Row 40 - Rennes
$ws.Range("U40").Value = 22

# Row 48 - VfL Wolfsburg
$ws.Range("U48").Value = 12

# Row 49 - Espanyol
$ws.Range("U49").Value = 12

# Row 51 - Osasuna
$ws.Range("U51").Value = 12

# Row 53 - Crystal Palace
$ws.Range("U53").Value = 2

# Row 54 - Brighton & Hove Albion
$ws.Range("U54").Value = 12

# Row 55 - TSG Hoffenheim
$ws.Range("U55").Value = 12

# Row 56 - Southampton
$ws.Range("U56").Value = 22

# Row 58 - Lille
$ws.Range("C58").Value = 10
$ws.Range("D58").Value = 10
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 3.57
$ws.Range("H58").Value = 0
$ws.Range("U58").Value = 8

# Row 62 - Bologna
$ws.Range("U62").Value = 12

# Row 65 - FC Union Berlin
$ws.Range("U65").Value = 2

# Row 68 - Monza
$ws.Range("U68").Value = 18

# Row 70 - FC Koln
$ws.Range("U70").Value = 2
